# Auto-generated Excel COM-interop edit script
# Applies the numeric corrections to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the commit "chore: update Sheets via scheduled runner".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1411.04
$ws.Range("J17").Value = 1411.04
$ws.Range("L17").Value = 4233.12
$ws.Range("N17").Value = -4569.12

$ws.Range("H98").Value = 4046.8667
$ws.Range("I98").Value = 1814.7142
$ws.Range("K98").Value = 1814.7142
$ws.Range("M98").Value = -316.7141999999999

$ws.Range("H106").Value = 2985.0715
$ws.Range("I106").Value = 3543.7778
$ws.Range("J106").Value = 1979.4
$ws.Range("K106").Value = 3543.7778
$ws.Range("L106").Value = 1979.4
$ws.Range("M106").Value = -2912.7778
$ws.Range("N106").Value = -3241.4

$ws.Range("H122").Value = 4046.8667
$ws.Range("I122").Value = 1814.7142
$ws.Range("K122").Value = 5444.142599999999
$ws.Range("M122").Value = -2994.142599999999

$ws.Range("H131").Value = 2881.389
$ws.Range("I131").Value = 2675.077
$ws.Range("J131").Value = 3417.8
$ws.Range("K131").Value = 8025.231000000001
$ws.Range("L131").Value = 10253.4
$ws.Range("M131").Value = -2985.231000000001
$ws.Range("N131").Value = -20333.4

$ws.Range("H132").Value = 1660.8422
$ws.Range("I132").Value = 1078.3438
$ws.Range("K132").Value = 3235.0314
$ws.Range("M132").Value = -705.0314000000003

$ws.Range("H137").Value = 4169.61
$ws.Range("I137").Value = 2616.1853
$ws.Range("K137").Value = 7848.5559
$ws.Range("M137").Value = -5298.5559

$ws.Range("H138").Value = 4581.695
$ws.Range("I138").Value = 3348.9092
$ws.Range("J138").Value = 5314.7026
$ws.Range("K138").Value = 10046.7276
$ws.Range("L138").Value = 15944.1078
$ws.Range("M138").Value = -4906.7276
$ws.Range("N138").Value = -26224.1078

$ws.Range("H141").Value = 6179.6665
$ws.Range("I141").Value = 4044
$ws.Range("K141").Value = 12132
$ws.Range("M141").Value = -6952

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5936.3335
$ws.Range("I2").Value = 5695.1577
$ws.Range("K2").Value = 5695.1577
$ws.Range("M2").Value = -5582.1577

$ws.Range("H45").Value = 7321.3
$ws.Range("I45").Value = 5040
$ws.Range("J45").Value = 8299
$ws.Range("K45").Value = 5040
$ws.Range("L45").Value = 8299
$ws.Range("M45").Value = -4663
$ws.Range("N45").Value = -9053

$ws.Range("H61").Value = 4956.3335
$ws.Range("I61").Value = 3722.2
$ws.Range("K61").Value = 3722.2
$ws.Range("M61").Value = -3510.2

$ws.Range("H62").Value = 45000
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46248

$ws.Range("H65").Value = 45000
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141240

$ws.Range("H74").Value = 1685.4073
$ws.Range("I74").Value = 1146.3334
$ws.Range("J74").Value = 5998
$ws.Range("K74").Value = 1146.3334
$ws.Range("L74").Value = 5998
$ws.Range("M74").Value = -272.3334
$ws.Range("N74").Value = -7746

$ws.Range("H77").Value = 1685.4073
$ws.Range("I77").Value = 1146.3334
$ws.Range("J77").Value = 5998
$ws.Range("K77").Value = 5731.666999999999
$ws.Range("L77").Value = 29990
$ws.Range("M77").Value = -1363.666999999999
$ws.Range("N77").Value = -38726

$ws.Range("H116").Value = 5936.3335
$ws.Range("I116").Value = 5695.1577
$ws.Range("K116").Value = 5695.1577
$ws.Range("M116").Value = -3401.1577

$ws.Range("H136").Value = 4956.3335
$ws.Range("I136").Value = 3722.2
$ws.Range("K136").Value = 11166.6
$ws.Range("M136").Value = -8616.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5936.3335
$ws.Range("I3").Value = 5695.1577
$ws.Range("K3").Value = 5695.1577
$ws.Range("M3").Value = -5581.1577

$ws.Range("H107").Value = 905.875
$ws.Range("I107").Value = 898.73334
$ws.Range("K107").Value = 898.73334
$ws.Range("M107").Value = 1021.26666

$ws.Range("H132").Value = 49999.57
$ws.Range("J132").Value = 49999.57
$ws.Range("L132").Value = 49999.57
$ws.Range("N132").Value = -60119.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 87075
$ws.Range("J74").Value = 104433.336
$ws.Range("L74").Value = 104433.336
$ws.Range("N74").Value = -106181.336

$ws.Range("H77").Value = 87075
$ws.Range("J77").Value = 104433.336
$ws.Range("L77").Value = 313300.008
$ws.Range("N77").Value = -322036.008

$ws.Range("H105").Value = 1715.5625
$ws.Range("I105").Value = 1629.9333
$ws.Range("K105").Value = 1629.9333
$ws.Range("M105").Value = 117.0667000000001

$ws.Range("H132").Value = 4471.72
$ws.Range("I132").Value = 3883.3333
$ws.Range("K132").Value = 11649.9999
$ws.Range("M132").Value = -9119.999899999999

$ws.Range("H134").Value = 3878.7576
$ws.Range("I134").Value = 3242.1365
$ws.Range("K134").Value = 9726.4095
$ws.Range("M134").Value = -7191.4095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 183.2
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 183.2
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 549.5999999999999
$ws.Range("N23").Value = -1019.6
$ws.Range("M23").ClearContents()

$ws.Range("H62").Value = 12149.5
$ws.Range("J62").Value = 12149.5
$ws.Range("L62").Value = 36448.5
$ws.Range("N62").Value = -37820.5

$ws.Range("H65").Value = 12149.5
$ws.Range("J65").Value = 12149.5
$ws.Range("L65").Value = 109345.5
$ws.Range("N65").Value = -116209.5

$ws.Range("H108").Value = 5583.091
$ws.Range("I108").Value = 5323.778
$ws.Range("K108").Value = 15971.334
$ws.Range("M108").Value = -13091.334

$ws.Range("H111").Value = 5131.75
$ws.Range("I111").Value = 4436.2856
$ws.Range("K111").Value = 13308.8568
$ws.Range("M111").Value = -10241.8568

$ws.Range("H112").Value = 3143.2
$ws.Range("I112").Value = 2554.25
$ws.Range("J112").Value = 5499
$ws.Range("K112").Value = 7662.75
$ws.Range("L112").Value = 16497
$ws.Range("M112").Value = -6554.75
$ws.Range("N112").Value = -18713

$ws.Range("H118").Value = 266
$ws.Range("I118").Value = 266
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 798
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 445
$ws.Range("N118").ClearContents()

$ws.Range("H119").Value = 1123.5
$ws.Range("I119").Value = 855.4286
$ws.Range("K119").Value = 2566.2858
$ws.Range("M119").Value = 2271.7142

$ws.Range("H133").Value = 7973
$ws.Range("I133").Value = 7973
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 23919
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -18859
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3256.3333
$ws.Range("I126").Value = 3064.25
$ws.Range("K126").Value = 9192.75
$ws.Range("M126").Value = -6722.75

$ws.Range("H135").Value = 112499.75
$ws.Range("J135").Value = 112499.75
$ws.Range("L135").Value = 112499.75
$ws.Range("N135").Value = -122639.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1970.5714
$ws.Range("I22").Value = 1466.8334
$ws.Range("K22").Value = 1466.8334
$ws.Range("M22").Value = -1171.8334

$ws.Range("H27").Value = 1970.5714
$ws.Range("I27").Value = 1466.8334
$ws.Range("K27").Value = 1466.8334
$ws.Range("M27").Value = -1359.8334

$ws.Range("H40").Value = 3004748
$ws.Range("I40").Value = 5003614.5
$ws.Range("J40").Value = 6448.625
$ws.Range("K40").Value = 5003614.5
$ws.Range("L40").Value = 6448.625
$ws.Range("M40").Value = -5003478.5
$ws.Range("N40").Value = -6720.625

$ws.Range("H132").Value = 6200
$ws.Range("I132").Value = 5500
$ws.Range("J132").Value = 6666.6665
$ws.Range("K132").Value = 16500
$ws.Range("L132").Value = 19999.9995
$ws.Range("M132").Value = -13970
$ws.Range("N132").Value = -25059.9995

$ws.Range("H136").Value = 6373.75
$ws.Range("I136").Value = 4750
$ws.Range("K136").Value = 14250
$ws.Range("M136").Value = -11700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 13063.75
$ws.Range("I4").Value = 12976.5
$ws.Range("J4").Value = 13500
$ws.Range("K4").Value = 12976.5
$ws.Range("L4").Value = 13500
$ws.Range("M4").Value = -12863.5
$ws.Range("N4").Value = -13726

$ws.Range("H8").Value = 3002
$ws.Range("J8").Value = 5004
$ws.Range("L8").Value = 5004
$ws.Range("N8").Value = -5284

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H119").Value = 60096.25
$ws.Range("J119").Value = 60096.25
$ws.Range("L119").Value = 60096.25
$ws.Range("N119").Value = -69772.25

$ws.Range("H132").Value = 3886.4666
$ws.Range("I132").Value = 2883.889
$ws.Range("J132").Value = 5390.3335
$ws.Range("K132").Value = 8651.667000000001
$ws.Range("L132").Value = 16171.0005
$ws.Range("M132").Value = -6121.667000000001
$ws.Range("N132").Value = -21231.0005

$ws.Range("H135").Value = 58069
$ws.Range("J135").Value = 58069
$ws.Range("L135").Value = 58069
$ws.Range("N135").Value = -68209
